$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text storage so values like
# "1.00" or "93.39" are not silently coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '44.872.75'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '2.255.75'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.73%  '
$ws.Range("D5").Value = '300.27'
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D6").Value = '93.39'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("E7").Value = '  -1.05%  '
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("D9").Value = '0.510'
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D10").Value = '33.95'
$ws.Range("E10").Value = '  -2.85%  '
$ws.Range("E11").Value = '  -2.33%  '
$ws.Range("D12").Value = '7.15'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").Value = '2.602.60'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '2.256.85'
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '13.62'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '0.793'
$ws.Range("E17").Value = '  -4.93%  '
$ws.Range("D18").Value = '44.742.89'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").Value = '12.66'
$ws.Range("E19").Value = '  +6.33%  '
$ws.Range("D20").Value = '0.0₃0917'
$ws.Range("E20").Value = '  -3.72%  '
$ws.Range("D21").Value = '6.05'
$ws.Range("E21").Value = '  -4.09%  '
$ws.Range("D22").Value = '65.14'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").Value = '239.34'
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  -4.91%  '
$ws.Range("B27").Value = 'InjectiveProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D27").Value = '39.47'
$ws.Range("E27").Value = '  +5.28%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '2.28'
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("D29").Value = '9.52'
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("D31").Value = '153.03'
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D32").Value = '5.55'
$ws.Range("E32").Value = '  -7.02%  '
$ws.Range("D33").Value = '0.0787'
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.117'
$ws.Range("E35").Value = '  -1.64%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.105'
$ws.Range("E36").Value = '  -3.60%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '2.88'
$ws.Range("E37").Value = '  -4.74%  '
$ws.Range("E38").Value = '  -6.28%  '
$ws.Range("D39").Value = '0.0300'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '3.70'
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Value = '3.23'
$ws.Range("E41").Value = '  -4.45%  '
$ws.Range("D42").Value = '13.80'
$ws.Range("E42").Value = '  -7.08%  '
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").Value = '1.786.86'
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").Value = '1.84'
$ws.Range("E45").Value = '  +8.94%  '
$ws.Range("D46").Value = '0.193'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").Value = '69.92'
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").Value = '75.49'
$ws.Range("E48").Value = '  -4.20%  '
$ws.Range("D49").Value = '95.99'
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("E50").Value = '  -4.64%  '
$ws.Range("D51").Value = '7.79'
$ws.Range("E51").Value = '  -2.81%  '
